$wb = $excel.ActiveWorkbook

# D (Options) and E (Total_Pnl) new values per sheet. B (Spot) becomes 0 for all rows; C (Fwds) stays 0.
$sheet2D = @(
    -608015.4219922652,
    -517084.9190198983,
    -433850.8881261189,
    -360371.4899629905,
    -297339.527727391,
    -243851.8482833975,
    -197819.4863662811,
    -156820.6600252705,
    -119028.6832030883,
    -83863.9511884173,
    -52175.86916700007,
    -25963.49722273639,
    -7801.011206403011,
    -194.3746579362669,
    -5061.328461508921,
    -23438.50359991641,
    -55425.98914322729,
    -100314.5082275465,
    -156813.9109880786,
    -223306.8687323929,
    -298074.018137013,
    -379462.9524502808,
    -465994.9563738405,
    -556417.024956425,
    -649712.8970714979,
    -745087.6603309007,
    -841938.2999796742,
    -939819.2437144688,
    -1038408.679978609,
    -1137478.773892827,
    -1236871.030688924
)
$sheet3D = @(
    -607083.6711008861,
    -513134.2711006302,
    -425976.0834276862,
    -348611.2604175693,
    -282895.109783117,
    -228726.0663800151,
    -184038.9964947673,
    -145624.876701955,
    -110397.4030100309,
    -76564.0796859269,
    -44259.8497388644,
    -15471.79102937753,
    6620.580257377544,
    18573.73569487785,
    17534.74941317584,
    1788.208730146667,
    -29076.39816945436,
    -74313.05541229132,
    -132351.4275373007,
    -201178.655513344,
    -278674.2058496354,
    -362850.05187899,
    -451988.8435148868,
    -544697.0677849774,
    -639899.5248764256,
    -736800.6195874116,
    -834832.1990431675,
    -933600.7666792932,
    -1032840.964115012,
    -1132377.986627994,
    -1232099.039241985
)
$sharedD = @(
    -612852.6396494271,
    -512852.6396494271,
    -412852.6396494271,
    -312852.6396494273,
    -212852.6396494273,
    -112852.6396494273,
    -12852.63964942723,
    87147.36035057277,
    -12852.63964942723,
    -112852.6396494272,
    -212852.6396494272,
    -112852.6396494272,
    -12852.6396494272,
    87147.3603505728,
    187147.3603505728,
    287147.3603505728,
    187147.3603505728,
    87147.3603505728,
    -12852.6396494272,
    -112852.6396494273,
    -212852.6396494274,
    -312852.6396494274,
    -412852.6396494274,
    -512852.6396494274,
    -612852.6396494271,
    -712852.639649427,
    -812852.639649427,
    -912852.639649427,
    -1012852.639649427,
    -1112852.639649427,
    -1212852.639649427
)

$ws = $wb.Worksheets.Item("today")
for ($i = 0; $i -lt $sheet2D.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = $sheet2D[$i]
    $ws.Cells.Item($row, 5).Value = $sheet2D[$i]
}

$ws = $wb.Worksheets.Item("EoM")
for ($i = 0; $i -lt $sheet3D.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = $sheet3D[$i]
    $ws.Cells.Item($row, 5).Value = $sheet3D[$i]
}

$ws = $wb.Worksheets.Item("2 Months")
for ($i = 0; $i -lt $sharedD.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = $sharedD[$i]
    $ws.Cells.Item($row, 5).Value = $sharedD[$i]
}

$ws = $wb.Worksheets.Item("3 Months")
for ($i = 0; $i -lt $sharedD.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = $sharedD[$i]
    $ws.Cells.Item($row, 5).Value = $sharedD[$i]
}

$ws = $wb.Worksheets.Item("4 Months")
for ($i = 0; $i -lt $sharedD.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = $sharedD[$i]
    $ws.Cells.Item($row, 5).Value = $sharedD[$i]
}

$ws = $wb.Worksheets.Item("5 Months")
for ($i = 0; $i -lt $sharedD.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = $sharedD[$i]
    $ws.Cells.Item($row, 5).Value = $sharedD[$i]
}

$ws = $wb.Worksheets.Item("6 Months")
for ($i = 0; $i -lt $sharedD.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = $sharedD[$i]
    $ws.Cells.Item($row, 5).Value = $sharedD[$i]
}

$ws = $wb.Worksheets.Item("Dec Expiry")
for ($i = 0; $i -lt $sharedD.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = $sharedD[$i]
    $ws.Cells.Item($row, 5).Value = $sharedD[$i]
}

